# Generate Report for Handback
#
# The localization run has come back: the status text moves from
# "Ready for handoff" to "Handed back: in sync with en-US" everywhere it
# appears (Overview + per-language sheets), and each per-language sheet's
# table gains two populated columns - F "Latest Target File" and
# G "Latest Handback File" - mirroring the source file name / handoff file
# name for each row. The "Latest Handback DateTime" column (H) is stamped
# with the real handback timestamp instead of the 0001-01-01 placeholder.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet: refresh the status text ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

# Row 2 (93748b00-0686-45dd-9908-0adf6aa0c8a2)
$wsZh.Range("F2").Value = $wsZh.Range("A2").Text
$wsZh.Range("F2").Style = "HyperLink"
$wsZh.Range("G2").Value = $wsZh.Range("D2").Text
$wsZh.Range("G2").Style = "HyperLink"

# Row 3 (cd2dc702-cf72-4b24-b41b-e33fc051edd5)
$wsZh.Range("F3").Value = $wsZh.Range("A3").Text
$wsZh.Range("F3").Style = "HyperLink"
$wsZh.Range("G3").Value = $wsZh.Range("D3").Text
$wsZh.Range("G3").Style = "HyperLink"

# Handback datetime actually recorded for the zh-cn run
$wsZh.Range("H2").Value = "2016-03-14 04:40:23"
$wsZh.Range("H3").Value = "2016-03-14 04:40:23"

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# Row 2 (93748b00-0686-45dd-9908-0adf6aa0c8a2)
$wsDe.Range("F2").Value = $wsDe.Range("A2").Text
$wsDe.Range("F2").Style = "HyperLink"
$wsDe.Range("G2").Value = $wsDe.Range("D2").Text
$wsDe.Range("G2").Style = "HyperLink"

# Row 3 (cd2dc702-cf72-4b24-b41b-e33fc051edd5)
$wsDe.Range("F3").Value = $wsDe.Range("A3").Text
$wsDe.Range("F3").Style = "HyperLink"
$wsDe.Range("G3").Value = $wsDe.Range("D3").Text
$wsDe.Range("G3").Style = "HyperLink"

# Handback datetime actually recorded for the de-de run
$wsDe.Range("H2").Value = "2016-03-14 04:40:29"
$wsDe.Range("H3").Value = "2016-03-14 04:40:29"
